$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns
$ws.Range("A1").Value = "prot_id"
$ws.Range("B1").Value = "new_prot_id"

# Swap order of TALA:TALB -> TALB:TALA in row 9
$ws.Range("A9").Value = "TALB:TALA"

# Update the active cell/selection to match the authored state
$ws.Range("A10").Select()
